# TIMBR_DAG.pptx revision: shrink the 8 callout labels around the DAG figure
# (Trait / Allele Effects / Diplotype States / Allelic Series / Tree /
# Branch Mutations / Concentration Parameter / Diplotype+Probabilities)
# to 10.5pt Arial and nudge their boxes to the re-flowed (autofit) size.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-LabelFont($shape) {
    $tr = $shape.TextFrame.TextRange
    $tr.Font.Size = 10.5
    $tr.Font.Name = "Arial"
    $tr.Font.NameComplexScript = "Arial"
}

function Set-LabelBox($shape, $left, $top, $width, $height) {
    $shape.Left = $left
    $shape.Top = $top
    $shape.Width = $width
    $shape.Height = $height
}

# 1) "Trait" textbox (TextBox 98 / id 99) -> Shapes.Item(61)
$shTrait = $s.Shapes.Item(61)
Set-LabelFont $shTrait
Set-LabelBox $shTrait 668.7217322834646 153.11779527559054 38.79283464566929 20.599212598425197

# 2) "Allele Effects" textbox (TextBox 99 / id 100) -> Shapes.Item(62)
$shAllele = $s.Shapes.Item(62)
Set-LabelFont $shAllele
Set-LabelBox $shAllele 129.86118320236218 72.64267716535433 74.83637795275591 32.71637925275591

# 3) "Diplotype States" textbox (TextBox 100 / id 101) -> Shapes.Item(63)
$shDiploStates = $s.Shapes.Item(63)
Set-LabelFont $shDiploStates
Set-LabelBox $shDiploStates 134.1048818897638 168.09543617086612 68.76755905511811 33.92811023622047

# 4) "Allelic Series" textbox (TextBox 101 / id 102) -> Shapes.Item(64)
$shAllelicSeries = $s.Shapes.Item(64)
Set-LabelFont $shAllelicSeries
Set-LabelBox $shAllelicSeries 128.9252755905512 265.14409448818895 79.1267738535433 33.92811023622047

# 5) "Tree" textbox (TextBox 102 / id 103) -> Shapes.Item(65)
$shTree = $s.Shapes.Item(65)
Set-LabelFont $shTree
Set-LabelBox $shTree 18.488976577952755 278.7762992125984 37.080393700787404 19.993386326771653

# 6) "Branch Mutations" textbox (TextBox 103 / id 104) -> Shapes.Item(66)
$shBranchMut = $s.Shapes.Item(66)
Set-LabelFont $shBranchMut
Set-LabelBox $shBranchMut 132.89425196850394 465.5155182110236 72.27622047244094 33.92811023622047

# 7) "Concentration Parameter" textbox (TextBox 104 / id 105) -> Shapes.Item(67)
$shConcParam = $s.Shapes.Item(67)
Set-LabelFont $shConcParam
Set-LabelBox $shConcParam 252.38858797716534 465.2146456692913 92.01149606299212 33.92811023622047

# 8) "Diplotype" / "Probabilities" two-paragraph textbox (TextBox 105 / id 106)
#    -> Shapes.Item(68). The font-name setters here only ever touch the
#    first run of the shape, so isolate the second paragraph as the sole
#    (first) paragraph, format it, then re-insert "Diplotype" + a hard
#    return in front of it -- the new leading run picks up the run
#    formatting of the paragraph it is inserted before.
$shDiploProb = $s.Shapes.Item(68)
$trDP = $shDiploProb.TextFrame.TextRange
$trDP.Text = "Probabilities"
Set-LabelFont $shDiploProb
[void]$trDP.InsertBefore("Diplotype`r")
Set-LabelBox $shDiploProb -1.7280314960629921 167.41464996929133 79.08464566929133 32.71637925275591
